# Reject appointment by interpreter is done.
# Add two new Login test-data rows (Interpreter / Invalid login) to the
# "Login" sheet, matching hyperlinks/styling of the existing rows, and move
# the active selection the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- Row 3: "Interpreter" --------------------------------------------------
$ws.Range("A3").Value2 = "Interpreter"
$ws.Range("B3").Value2 = "wei.yuan@sstech.us"
$ws.Range("C3").Value2 = "Welcome@1"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:wei.yuan@sstech.us", [Type]::Missing, [Type]::Missing, "wei.yuan@sstech.us") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Welcome@1", [Type]::Missing, [Type]::Missing, "Welcome@1") | Out-Null

# --- Row 4: "Invalid login" (plain text, no hyperlinks) --------------------
$ws.Range("A4").Value2 = "Invalid login"
$ws.Range("B4").Value2 = "abc.sstech.us"
$ws.Range("C4").Value2 = "xyz"

# --- Formatting: keep the same look as the rest of the table ---------------
# Header row stays bold.
$ws.Range("A1:C1").Font.Bold = $true

# Column A (Description) data cells are plain, unlinked text.
$ws.Range("A2:A4").Font.Bold = $false
$ws.Range("A2:A4").Font.Underline = $false

# Columns B/C (EmailAddress/Password) look like the hyperlink cells above
# them, underlined, even for the rows that are not real mailto links.
$ws.Range("B4:C4").Font.Underline = $true
$ws.Range("B5:B6").Font.Underline = $true

# --- Leave the selection the way it was when the workbook was last saved ---
$ws.Range("E10").Select() | Out-Null
